# Editing Melanie's sulfide code to read in multiple plates.
#
# A new plate ("Plate 3 20251114") was inserted ahead of the existing
# 20251114 plates, so the old "Plate 3/4/5 20251114" tabs each shift up
# by one name: Plate 3 -> Plate 4, Plate 4 -> Plate 5, Plate 5 -> Plate 6.
# Rename starting from the last sheet so we never collide with a name
# that hasn't been vacated yet.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Plate 5 20251114").Name = "Plate 6 20251114"
$wb.Worksheets.Item("Plate 4 20251114").Name = "Plate 5 20251114"
$wb.Worksheets.Item("Plate 3 20251114").Name = "Plate 4 20251114"

# The newly-renamed "Plate 4 20251114" tab (the former "Plate 3 20251114")
# is now the active/selected tab, with cell S24 selected.
$ws = $wb.Worksheets.Item("Plate 4 20251114")
$ws.Activate() | Out-Null
$ws.Range("S24").Select() | Out-Null
